$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "70.774.48"
$c.ClearFormats()
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.20%  "
$c.ClearFormats()

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.532.68"
$c.ClearFormats()
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.39%  "
$c.ClearFormats()

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.15%  "
$c.ClearFormats()

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "622.61"
$c.ClearFormats()
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.74%  "
$c.ClearFormats()

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "174.19"
$c.ClearFormats()
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.53%  "
$c.ClearFormats()

$ws.Cells.Item(7, 2).Value = "XRP"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.613"
$c.ClearFormats()
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.89%  "
$c.ClearFormats()

$ws.Cells.Item(8, 2).Value = "LidoStakedEther"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "3.530.12"
$c.ClearFormats()
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.39%  "
$c.ClearFormats()

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c.ClearFormats()

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.199"
$c.ClearFormats()
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.49%  "
$c.ClearFormats()

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "7.07"
$c.ClearFormats()
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.76%  "
$c.ClearFormats()

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.584"
$c.ClearFormats()
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.33%  "
$c.ClearFormats()

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "46.54"
$c.ClearFormats()
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.35%  "
$c.ClearFormats()

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000278"
$c.ClearFormats()
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.49%  "
$c.ClearFormats()

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "4.096.72"
$c.ClearFormats()
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.65%  "
$c.ClearFormats()

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "8.41"
$c.ClearFormats()
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.59%  "
$c.ClearFormats()

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "609.88"
$c.ClearFormats()
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.73%  "
$c.ClearFormats()

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "3.534.38"
$c.ClearFormats()
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.18%  "
$c.ClearFormats()

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "70.840.04"
$c.ClearFormats()
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.22%  "
$c.ClearFormats()

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "17.77"
$c.ClearFormats()
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.38%  "
$c.ClearFormats()

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.77%  "
$c.ClearFormats()

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "9.10"
$c.ClearFormats()
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.52%  "
$c.ClearFormats()

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "98.60"
$c.ClearFormats()
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.70%  "
$c.ClearFormats()

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "15.61"
$c.ClearFormats()
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.82%  "
$c.ClearFormats()

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.25%  "
$c.ClearFormats()

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.09%  "
$c.ClearFormats()

$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.13%  "
$c.ClearFormats()

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "33.79"
$c.ClearFormats()
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.43%  "
$c.ClearFormats()

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "9.09"
$c.ClearFormats()
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.22%  "
$c.ClearFormats()

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "3.02"
$c.ClearFormats()
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.89%  "
$c.ClearFormats()

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "8.12"
$c.ClearFormats()
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.91%  "
$c.ClearFormats()

$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.53%  "
$c.ClearFormats()

$ws.Cells.Item(34, 2).Value = "Bittensor"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "634.04"
$c.ClearFormats()
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.20%  "
$c.ClearFormats()

$ws.Cells.Item(35, 2).Value = "NEARProtocol"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "6.81"
$c.ClearFormats()
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.07%  "
$c.ClearFormats()

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.0997"
$c.ClearFormats()
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.82%  "
$c.ClearFormats()

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "10.81"
$c.ClearFormats()
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.93%  "
$c.ClearFormats()

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.0475"
$c.ClearFormats()
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.04%  "
$c.ClearFormats()

$ws.Cells.Item(39, 2).Value = "dogwifhat"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "3.44"
$c.ClearFormats()
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  -9.76%  "
$c.ClearFormats()

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "57.04"
$c.ClearFormats()
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.76%  "
$c.ClearFormats()

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.38%  "
$c.ClearFormats()

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.144"
$c.ClearFormats()
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.58%  "
$c.ClearFormats()

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "3.354.77"
$c.ClearFormats()
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.73%  "
$c.ClearFormats()

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0732"
$c.ClearFormats()
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.83%  "
$c.ClearFormats()

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.56%  "
$c.ClearFormats()

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.46%  "
$c.ClearFormats()

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "32.00"
$c.ClearFormats()
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.77%  "
$c.ClearFormats()

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.62%  "
$c.ClearFormats()

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.82%  "
$c.ClearFormats()

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "133.09"
$c.ClearFormats()
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.ClearFormats()

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.ClearFormats()
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +5.19%  "
$c.ClearFormats()
